# Update the "Förändrad" (Changed) date column (C) for rows 2..205
# from 2023-09-17 (serial 45186) to 2023-09-19 (serial 45188).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 205; $row++) {
    $ws.Cells.Item($row, 3).Value = 45188
}
